# Generate Report for Handoff
# The a9227bf9-d672-45db-85a6-1fe97592d078 file is now ready for handoff;
# update its status + handoff timestamps across the Overview and per-locale
# (zh-cn / de-de) sheets.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"

# --- Overview sheet: row 3 is the a9227bf9-d672-...-...md entry ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status
$overview.Range("D3").Value = "2016-38-20 04:38:15"

# --- zh-cn sheet: row 3 is the a9227bf9-d672-...-...md entry ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $status
$zhcn.Range("E3").Value = "2016-03-20 04:38:12"

# --- de-de sheet: row 3 is the a9227bf9-d672-...-...md entry ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $status
$dede.Range("E3").Value = "2016-03-20 04:38:15"
